$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# repull data, push all data, mean calculation
# Update column F (dSF) values that changed after repulling/pushing data
$ws.Range("F2").Value = 3
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = -2
$ws.Range("F6").Value = 3
$ws.Range("F7").Value = -5
$ws.Range("F8").Value = -1
$ws.Range("F14").Value = -3
$ws.Range("F15").Value = 0
$ws.Range("F18").Value = 5
